# Created method for generating tiles for layer 1, untested
#
# Updates the "Vera" memory allocation worksheet:
#  - Map Base segment grows from 2561 to 4096 bytes
#  - The old "Vacant" gap row (row 7) is repurposed into a new
#    "Sprites" allocation entry (4097 bytes)
#  - The previous "Sprites" (row 8) and "Volatile Buffer" (row 9)
#    rows are removed entirely, shifting the summary rows up
#  - The active selection moves to F8

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Vera")

# Grow the "Map Base" allocation size (row 6)
$ws.Range("C6").Value = 4096

# Remove the old "Sprites" row (8) and "Volatile Buffer" row (9).
# Deleting shifts the summary rows (Total Used / Total Available /
# Remaining) up automatically and keeps their formulas correct.
$ws.Rows.Item(8).EntireRow.Delete() | Out-Null
$ws.Rows.Item(8).EntireRow.Delete() | Out-Null

# Turn the former "Vacant" row (7) into the new "Sprites" entry.
# Its Start/Address/Ends formulas already reference the row above,
# so only the label and size need to change.
$ws.Range("A7").Value = "Sprites"
$ws.Range("C7").Value = 4097

# Match the author's final selection
$ws.Range("F8").Select() | Out-Null
